$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Email, CNP, Telefon) in F1:H1
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "CNP"
$ws.Range("H1").Value = "Telefon"

# Give the new headers the same formatting as the existing header row
# (copy format only, so the new cells reuse the existing bold header style)
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Resize the first, second and fifth columns to fit the new, narrower layout
$ws.Columns.Item(1).ColumnWidth = 26.92
$ws.Columns.Item(2).ColumnWidth = 32.64
$ws.Columns.Item(5).ColumnWidth = 11.6

# Move the active selection to the newly added H1 header cell
$ws.Range("H1").Select()
